$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Update the metadata values (URL, Identifier, Version) ---
$ws.Range("B2").Value = "http://fhir.ch/ig/ch-epr-term/ValueSet/DocumentEntry.typeCode"
$ws.Range("B3").Value = "id: 2.16.756.5.30.1.127.3.10.1.27 (use: OFFICIAL)"
$ws.Range("B4").Value = "2.0.0-ballot"

# --- Insert a new "Contact" row after the existing one (row 11), pushing
#     Jurisdiction/Description/Purpose/Copyright/Immutable down by one ---
$ws.Rows.Item(12).Insert()

# Copy formatting (borders/fill/font/alignment) from row 11 onto the new row 12
$ws.Range("A11:B11").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122)

# Populate the new row's values (duplicate of the Contact row above it)
$ws.Range("A12").Value = "Contact"
$ws.Range("B12").Value = "No display for ContactDetail"
